$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, preserving the "no explicit style" state
# that the source inlineStr cells have (avoids numeric auto-coercion of strings like
# "6.50" / "0.999" / "325.70" while not leaving a stray style index behind).
$refStyle = $ws.Range("C2").Style
function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $refStyle
}

Set-TextValue "D2" '66.744.27'
Set-TextValue "D3" '3.519.15'
Set-TextValue "E3" '  +0.58%  '
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '584.83'
Set-TextValue "E5" '  -2.48%  '
Set-TextValue "D6" '176.13'
Set-TextValue "E6" '  -2.44%  '
Set-TextValue "D8" '3.515.53'
Set-TextValue "E8" '  +0.42%  '
Set-TextValue "E9" '  -2.45%  '
Set-TextValue "E10" '  -2.72%  '
Set-TextValue "D11" '6.92'
Set-TextValue "E11" '  -1.92%  '
Set-TextValue "E12" '  -3.12%  '
Set-TextValue "D13" '4.124.26'
Set-TextValue "E13" '  +0.57%  '
Set-TextValue "D14" '30.58'
Set-TextValue "E14" '  -5.38%  '
Set-TextValue "E15" '  -1.58%  '
Set-TextValue "D16" '66.712.64'
Set-TextValue "E16" '  -1.16%  '
Set-TextValue "E17" '  -2.56%  '
Set-TextValue "D18" '3.510.34'
Set-TextValue "E18" '  +0.40%  '
Set-TextValue "E19" '  -3.93%  '
Set-TextValue "D20" '14.02'
Set-TextValue "E20" '  -2.03%  '
Set-TextValue "D21" '382.47'
Set-TextValue "E21" '  -2.07%  '
Set-TextValue "D22" '7.92'
Set-TextValue "E22" '  -0.60%  '
Set-TextValue "D23" '0.552'
Set-TextValue "E23" '  +1.69%  '
Set-TextValue "E24" '  +0.12%  '
Set-TextValue "D25" '72.38'
Set-TextValue "E25" '  -1.98%  '
Set-TextValue "E26" '  +0.29%  '
Set-TextValue "E27" '  -1.13%  '
Set-TextValue "D28" '9.91'
Set-TextValue "E28" '  -4.72%  '
Set-TextValue "E29" '  -1.62%  '
Set-TextValue "D30" '0.999'
Set-TextValue "E30" '  -0.02%  '
Set-TextValue "D31" '24.63'
Set-TextValue "E31" '  +4.35%  '
Set-TextValue "D32" '5.94'
Set-TextValue "E32" '  -4.05%  '
Set-TextValue "E33" '  -2.80%  '
Set-TextValue "E34" '  -5.43%  '
Set-TextValue "D35" '7.28'
Set-TextValue "E35" '  -1.97%  '
Set-TextValue "D36" '0.999'
Set-TextValue "E36" '  -0.05%  '
Set-TextValue "E37" '  -2.10%  '
Set-TextValue "D38" '30.18'
Set-TextValue "E38" '  +13.73%  '
Set-TextValue "D39" '161.47'
Set-TextValue "E39" '  -0.77%  '
Set-TextValue "E40" '  +2.44%  '
Set-TextValue "E41" '  -5.16%  '
Set-TextValue "D42" '4.54'
Set-TextValue "E42" '  -2.47%  '
Set-TextValue "D43" '6.50'
Set-TextValue "E43" '  -5.09%  '
Set-TextValue "D44" '2.56'
Set-TextValue "E44" '  -9.60%  '
Set-TextValue "D45" '2.730.05'
Set-TextValue "E45" '  -4.49%  '
Set-TextValue "D46" '0.0704'
Set-TextValue "E46" '  -2.82%  '
Set-TextValue "D47" '40.76'
Set-TextValue "D48" '25.12'
Set-TextValue "E48" '  -6.52%  '
Set-TextValue "E49" '  -2.65%  '
Set-TextValue "D50" '325.70'
Set-TextValue "E50" '  -2.44%  '
Set-TextValue "E51" '  -3.73%  '
